$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.87893676757812
$ws.Range("C3").Value = 17.11583137512207
$ws.Range("C4").Value = 16.7238712310791
$ws.Range("C5").Value = 17.0588493347168
$ws.Range("C6").Value = 16.89910888671875
